$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "luận giải" (interpretation) rows appended below the existing table.
# Column A values are entered first (top to bottom), then the column B
# values are filled back in (matching the order the author originally typed
# them), so new shared-string entries land in the same sequence as the
# source workbook.

$ws.Range("A11").Value = "Tử Phủ Vũ Tướng"
$ws.Range("A12").Value = "Kình Đà Không Kiếp"
$ws.Range("B12").Value = "Đại vận rất xấu, sẽ gặp nhiều điều khó khăn, trở ngại "
$ws.Range("A13").Value = "Phá Không Kiếp"
$ws.Range("A14").Value = "Sát Kình Đà"
$ws.Range("A15").Value = "Tham Hỏa Linh"
$ws.Range("B14").Value = "Mệnh có Thất Sát, hội hợp chiếu với Kình Dương, Đà La : Đại vận xấu tuy nhiên có phần giảm nhẹ do Thất Sát thủ mệnh."
$ws.Range("B13").Value = "Mệnh có Phá Quân, hội hợp chiếu với Địa Không, Địa Kiếp: Đại vận cực xấu tuy nhiên cũng có phần giảm nhẹ do Phá Quân thủ mệnh."
$ws.Range("B15").Value = "Mệnh có Tham Lang, hội hợp chiếu với Hỏa Tinh, Linh Tinh: Đại vận xấu tuy nhiên có phần giảm nhẹ do Tham Lang thủ mệnh."
$ws.Range("B11").Value = "Đại vận cần bản lĩnh và kĩ năng lãnh đạo."
$ws.Range("A16").Value = "Cơ Nguyệt Đồng Lương"
$ws.Range("A17").Value = "Cự Nhật"
$ws.Range("B17").Value = "Đại vận có nhiều thị phi qua lời nói, dù là đấu tranh chỉ để tìm đến điều đúng đắn nhưng cũng rất vất vả."
$ws.Range("A18").Value = "Tướng Triệt Hình Kỵ"
$ws.Range("B18").Value = "Đại vận xấu về đường quan lộc."
$ws.Range("A19").Value = "Lưu Hà Kiếp Sát"
$ws.Range("B19").Value = "Đại vận rất xấu"

# Match the author's final selection (cell below/right of the new table).
$ws.Range("B22").Select() | Out-Null
